# "Sluiten" toegevoegd als bevestiging van projectleiders dat iets gesloten
# kan worden: insert a new "Sluiten" column between "Algemene informatie"
# and "Actiepunten Bram" on the main sheet's header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at J, shifting the existing J:N columns to K:O.
# Excel's default insert behaviour (format-from-left-neighbour) reproduces
# the exact style ids seen in the target file (new J1 -> style of I1,
# shifted former J1 -> K1 keeps its own style).
$ws.Columns("J").Insert() | Out-Null

# Give the new column its header text in the second (header-label) row.
$ws.Range("J2").Value = "Sluiten"

# Match the saved selection/active cell from the edit.
$ws.Range("J3").Select() | Out-Null
